$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "CasesTab" row label is actually the Participants tab - fix it.
$ws.Range("A2").Value = "ParticipantsTab"

# Move the active selection to A2 (top of the renamed row).
$ws.Activate()
$ws.Range("A2").Select()
